$d = $word.ActiveDocument

# --- Step 1: replace paragraph 1 (index 1) title text ---
$d.Paragraphs(1).Range.Text = "1) European Parliament: Trends of Political Groups and"

# --- Step 2: insert new paragraph 'The main commissions' after paragraph 1 ---
$d.Paragraphs(1).Range.InsertParagraphAfter()
$d.Paragraphs(2).Range.Text = "The main commissions"

# --- Step 3: paragraph 3 ("Introduction", originally paragraph 2) is unchanged ---

# --- Step 4: replace paragraphs (originally 2..22, now 4..24) text 1:1 in order ---
$d.Paragraphs(4).Range.Text = "In the context of the last elections to the European Parliament, which were held on June 8 and 9, we will return in this article on the main issues of existing political groups as well as the major committees present."
$d.Paragraphs(5).Range.Text = "In fact, the stakes concerning these elections seem numerous, with many issues raised, such as the environment, the economy, the still present consequences of the pandemic, immigration, national sovereignty or even the field of digital and new technologies, everyone's rights, equality and finally the CAP (Common Agricultural Policy)."
$d.Paragraphs(6).Range.Text = "The European Parliament has a number of functions to date, crucial for the proper functioning of Europe as a whole. The functions are both legislative and democratic, but also budgetary and relational, with the desire to highlight a number of values such as human, solidarity and equality."
$d.Paragraphs(7).Range.Text = "The trends of the different European political groups."
$d.Paragraphs(8).Range.Text = "There are to date eight main groups in the European Parliament."
$d.Paragraphs(9).Range.Text = "The European People's Party, EPP"
$d.Paragraphs(10).Range.Text = "- The Renew Europe, RE"
$d.Paragraphs(11).Range.Text = "The Identity and Democracy Group, ID"
$d.Paragraphs(12).Range.Text = "Progressive Alliance of Socialists and Democrats"
$d.Paragraphs(13).Range.Text = "The Left Group in the European Parliament and the Nordic Green Left"
$d.Paragraphs(14).Range.Text = "- European Conservative Reformists, ECR"
$d.Paragraphs(15).Range.Text = "The Greens and the European Free Alliance, EFA"
$d.Paragraphs(16).Range.Text = "The unregistered"
$d.Paragraphs(17).Range.Text = "The European People's Party, or EPP"
$d.Paragraphs(18).Range.Text = "It is a right-wing party / center-right, which includes conservatives and essentially Christian Democrats."
$d.Paragraphs(19).Range.Text = "They constitute a majority in the European Parliament, it is a very important group. In the last elections, the group represents 189 seats and 26.25% of the votes. They are mostly in 12 Member States such as Spain but also Germany and Poland. The PPE attaches primary importance to immigration control as well as aid to farmers. Security is also a priority."
$d.Paragraphs(20).Range.Text = "Renew Europe, or RE"
$d.Paragraphs(21).Range.Text = "This group belongs to the center/liberal."
$d.Paragraphs(22).Range.Text = "It is an alliance between several groups, the Democrats and Liberals as well as the Republic on the Move, founded by Emmanuel Macron."
$d.Paragraphs(23).Range.Text = "The group records in 2024 a number of 79 seats for 10.79% of the votes."
$d.Paragraphs(24).Range.Text = "The main trends are unity and collaboration as well as European security. The group also emphasizes energy transition, and the fight against extremism, with a refusal to develop an alliance with the far right."

# --- Step 5: append new paragraphs at the end of the document ---
$lastIdx = $d.Paragraphs.Count
$d.Paragraphs($lastIdx).Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$d.Paragraphs($lastIdx).Range.Text = "Identity and Democracy, or ID"
$d.Paragraphs($lastIdx).Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$d.Paragraphs($lastIdx).Range.Text = "It is a right-wing group, a link is established with the Italian League party and the National Rally. The results are quite significant in 2024 since they have 58 seats and 8.06% of the vote at the last election."
$d.Paragraphs($lastIdx).Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$d.Paragraphs($lastIdx).Range.Text = "The group attaches importance to work and safety and clearly positions itself against immigration."
$d.Paragraphs($lastIdx).Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$d.Paragraphs($lastIdx).Range.Text = "Progressive Alliance of Socialists and Democrats (S&D)"
$d.Paragraphs($lastIdx).Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$d.Paragraphs($lastIdx).Range.Text = "At the last elections, the group has 135 seats and 18.75% of votes."
$d.Paragraphs($lastIdx).Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$d.Paragraphs($lastIdx).Range.Text = "It is the second largest group after the EPP."
$d.Paragraphs($lastIdx).Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$d.Paragraphs($lastIdx).Range.Text = "1"
$d.Paragraphs($lastIdx).Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$d.Paragraphs($lastIdx).Range.Text = "European United Left/Nordic Green Left"
$d.Paragraphs($lastIdx).Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$d.Paragraphs($lastIdx).Range.Text = "The group accounts for in 2024 a number of 36 seats and 5% of the votes."
$d.Paragraphs($lastIdx).Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$d.Paragraphs($lastIdx).Range.Text = "They advocate social justice, environmental protection and also the assiduous respect of human rights in all countries of Europe. They are for feminism and equality throughout the territory."
$d.Paragraphs($lastIdx).Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$d.Paragraphs($lastIdx).Range.Text = "European Conservatives and Reformists Group"
$d.Paragraphs($lastIdx).Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$d.Paragraphs($lastIdx).Range.Text = "With 73 seats and 10.14% of the vote in 2024, the group includes Vox in Spain and the Democrats of Sweden among its ranks."
$d.Paragraphs($lastIdx).Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$d.Paragraphs($lastIdx).Range.Text = "They are for a perfect equality between the Member States, and wish to find more radical and more effective solutions against terrorism and all forms of crime. They are for national sovereignty and for an overall more responsible Europe."
$d.Paragraphs($lastIdx).Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$d.Paragraphs($lastIdx).Range.Text = "The Greens / Free European Alliance"
$d.Paragraphs($lastIdx).Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$d.Paragraphs($lastIdx).Range.Text = "The Greens get 53 seats in the last elections and 7,36% of votes. They want to fight more virulently against climate change and against threats related to sexism, gender especially. They are also in favor of social justice and solidarity."
$d.Paragraphs($lastIdx).Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$d.Paragraphs($lastIdx).Range.Text = "The unregistered"
$d.Paragraphs($lastIdx).Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$d.Paragraphs($lastIdx).Range.Text = "They are the deputies who do not wish to belong to any of the groups mentioned above. They get 45 seats and 6.25% of the vote. These include former members of the National Rally as well as Catalan separatists. They have both independence in voting and decisions made."
$d.Paragraphs($lastIdx).Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$d.Paragraphs($lastIdx).Range.Text = "The main committees of the European Parliament"
$d.Paragraphs($lastIdx).Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$d.Paragraphs($lastIdx).Range.Text = "Here are some present commissions in Parliament, the list is not exhaustive."
$d.Paragraphs($lastIdx).Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$d.Paragraphs($lastIdx).Range.Text = "1"
$d.Paragraphs($lastIdx).Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$d.Paragraphs($lastIdx).Range.Text = "1"
$d.Paragraphs($lastIdx).Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$d.Paragraphs($lastIdx).Range.Text = "1) The INTA is a specialized committee in charge of international trade. It deals with trade agreements between countries. 2) The INTA, a specialised committee on international trade. It deals with trade agreements between countries."
$d.Paragraphs($lastIdx).Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$d.Paragraphs($lastIdx).Range.Text = "The BUDG is the budget committee, which deals with financial matters throughout the EU."
$d.Paragraphs($lastIdx).Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$d.Paragraphs($lastIdx).Range.Text = "The EMPL Committee is the one relating to employment and social affairs. It takes charge of all policies related to employment as well as health."
$d.Paragraphs($lastIdx).Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$d.Paragraphs($lastIdx).Range.Text = "The ENVI Committee is the committee on environment and food safety."
$d.Paragraphs($lastIdx).Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$d.Paragraphs($lastIdx).Range.Text = "The IMCO is the committee on consumer affairs and their protection as well as more generally to the internal market. Given the extent of the online market, these committees take a considerable place in economic issues."
$d.Paragraphs($lastIdx).Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$d.Paragraphs($lastIdx).Range.Text = "AGRI is the committee on agriculture and rural development. These are very topical issues that interest consumers in all Member States more and more."
$d.Paragraphs($lastIdx).Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$d.Paragraphs($lastIdx).Range.Text = "1"
$d.Paragraphs($lastIdx).Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$d.Paragraphs($lastIdx).Range.Text = "There are other ones, such as TRAN, which is the committee in charge of tourism and the transport sector or even the REGI committee, which deals more with everything that concerns regions."
$d.Paragraphs($lastIdx).Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$d.Paragraphs($lastIdx).Range.Text = "1) CONCLUSION"
$d.Paragraphs($lastIdx).Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$d.Paragraphs($lastIdx).Range.Text = "The European Parliament represents a major importance in the whole of European legislation. It deals with citizens of Europe, their interests but also their well-being, their rights as well as their duties. The Parliament is composed of several groups and commissions which have different ideas, and whose objectives differ according to political ideals. The objective of the Parliament as a whole is to respond more effectively to all Europeans, to various problems such as transport, tourism, environment or security."
$d.Paragraphs($lastIdx).Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$d.Paragraphs($lastIdx).Range.Text = "Le DEVE, commission du développement, traite de l’humanitaire et de la coopération entre États."
$d.Paragraphs($lastIdx).Range.InsertParagraphAfter()
$lastIdx = $lastIdx + 1
$d.Paragraphs($lastIdx).Range.Text = "L’AFET, commission des affaires étrangères, qui s’occupe comme son nom l’indique des relations avec les autres pays membres de l’UE, surtout concernant la sécurité et la défense nationale."
